# "Matrix Rain Effect added to the report module"
# Refresh the sample PO_Detail row with a new Quote/DA/Unit combo and clear
# out the stale "In Service Date" value, then leave PO_Detail as the
# selected tab (instead of Unit_to_Reconcile_Output).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_Detail")

# C2 is about to take over D2's old "text / quote-prefixed" number format
# (D2 keeps that format but its value goes away below), so copy the format
# across first.
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Stage each new value as text on a scratch cell and paste-special just the
# value into place, so the target cell's existing number format/style is
# left alone (A2, C2) or no style is forced at all (B2).
$scratch = $ws.Range("Z100")

$scratch.Value = "'382425"
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$scratch.Value = "'1"
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$scratch.Value = "'00999166"
$scratch.Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

# The old "In Service Date" (07/01/2017) is simply removed; the cell keeps
# whatever formatting it had.
$ws.Range("D2").ClearContents()

# Make PO_Detail the active/selected sheet.
$ws.Activate()
